$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I5").Value = -0.6575302404215204
$ws.Range("J5").Value = 0.4508445361474317
$ws.Range("K5").Value = 0.09288533878809677
$ws.Range("L5").Value = 2.496184122438338
